$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.417.58'
$ws.Range('E2').Value = '  +2.62%  '
$ws.Range('D3').Value = '3.187.21'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '596.34'
$ws.Range('E5').Value = '  +3.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.23'
$ws.Range('E6').Value = '  +3.90%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.183.34'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.549'
$ws.Range('E9').Value = '  +4.65%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.95'
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.517'
$ws.Range('E12').Value = '  +3.85%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000269'
$ws.Range('E13').Value = '  +3.49%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '39.19'
$ws.Range('E14').Value = '  +5.81%  '
$ws.Range('D15').Value = '3.704.76'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '66.345.77'
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.47'
$ws.Range('E17').Value = '  +5.17%  '
$ws.Range('D18').Value = '3.185.01'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '514.90'
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.38'
$ws.Range('E21').Value = '  +3.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.738'
$ws.Range('E22').Value = '  +3.63%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.06'
$ws.Range('E23').Value = '  +4.79%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '14.94'
$ws.Range('E24').Value = '  -2.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '85.65'
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.28'
$ws.Range('E27').Value = '  +4.62%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.00'
$ws.Range('E28').Value = '  +3.87%  '
$ws.Range('E29').Value = '  +7.66%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.13'
$ws.Range('E30').Value = '  +15.97%  '
$ws.Range('E31').Value = '  +4.62%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '28.30'
$ws.Range('E32').Value = '  +3.18%  '
$ws.Range('E33').Value = '  +2.67%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.53'
$ws.Range('E35').Value = '  +1.39%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '507.75'
$ws.Range('E36').Value = '  +6.80%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '54.88'
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0901'
$ws.Range('E38').Value = '  +0.99%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0425'
$ws.Range('E39').Value = '  +2.73%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.127'
$ws.Range('E40').Value = '  +10.06%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '8.91'
$ws.Range('E41').Value = '  +2.87%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.304'
$ws.Range('E42').Value = '  +8.47%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.88'
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').Value = '0.0₃0676'
$ws.Range('E44').Value = '  +16.74%  '
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('D46').Value = '2.909.47'
$ws.Range('E46').Value = '  -3.03%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '28.67'
$ws.Range('E47').Value = '  +2.30%  '
$ws.Range('E48').Value = '  +3.17%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.34'
$ws.Range('E50').Value = '  +5.30%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.64'
$ws.Range('E51').Value = '  +9.40%  '
